$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting existing rows 126-175 down to 127-176.
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new weekly record.
$ws.Range("A126").Value = 9
$ws.Range("B126").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C126").Value = "Metropolitana"
$ws.Range("D126").Value = 44489
$ws.Range("E126").Value = 13
$ws.Range("F126").Value = 100112021
$ws.Range("G126").Value = "Ají"
$ws.Range("H126").Value = "Americana (o)"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 7
$ws.Range("K126").Value = 50000
$ws.Range("L126").Value = 52000
$ws.Range("M126").Value = 50857
$ws.Range("N126").Value = '$/caja 25 kilos'
$ws.Range("O126").Value = "Provincia de Limarí"
$ws.Range("P126").Value = 2034
$ws.Range("Q126").Value = 25
$ws.Range("R126").Value = "Hortaliza"
